# #updates in GUI and Excel settings
#
# Rename "Test Cases" -> "TestCases" and "Test Steps" -> "TestSteps",
# fix a couple of header labels / a typo'd keyword name, resize a
# column on the Test Steps sheet, and update the remembered
# selection/active-tab state for each sheet (as last left by the user).

$wb = $excel.ActiveWorkbook

# --- Rename sheets -----------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Name = "TestCases"

$wsTestSteps = $wb.Worksheets.Item("Test Steps")
$wsTestSteps.Name = "TestSteps"

# --- Settings sheet: fix keyword typo -----------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("A16").Value = "AssertTextPresent"

# --- TestCases sheet: header rename -------------------------------------
$wsTestCases.Range("E1").Value = "Result"

# --- TestSteps sheet: header renames + keyword typo fix + column width --
$wsTestSteps.Range("A1").Value = "TestCaseID"
$wsTestSteps.Range("B1").Value = "StepNo"
$wsTestSteps.Range("D1").Value = "PageName"
$wsTestSteps.Range("E1").Value = "PageObject"
$wsTestSteps.Range("H1").Value = "Result"
$wsTestSteps.Range("F6").Value = "AssertTextPresent"
$wsTestSteps.Range("F14").Value = "AssertTextPresent"
$wsTestSteps.Columns.Item(7).ColumnWidth = 22

# --- Restore each sheet's last-used selection (also drives which sheet
#     ends up active/tabSelected + the workbook's activeTab) ------------
$wsSettings.Range("C18").Select()
$wb.Worksheets.Item("Pages Objects").Range("C5").Select()
$wb.Worksheets.Item("Locators").Range("F20").Select()
$wsTestCases.Range("F11").Select()
$wsTestSteps.Range("F14").Select()
